$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: rows 6 and 7 (source files 34fcff7b..., 7ddb50d5...) just got a new
# handoff cycle: Priority is now "ht" and the Latest Handoff Datetime refreshed.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = "2017-03-02 08:29:51"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H7").Value = "2017-03-02 08:29:51"

# --- de-de sheet: rows 4,5,6,7 (source files 167f43ef..., 34fcff7b..., 35cc0030..., 7ddb50d5...)
# also get the "ht" priority marker, and their Latest Handoff Datetime (shared with the
# Overview sheet's generation timestamp) moves forward too.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("H4").Value = "2017-03-02 08:30:09"
$wsDeDe.Range("H5").Value = "2017-03-02 08:30:09"
$wsDeDe.Range("H6").Value = "2017-03-02 08:30:09"
$wsDeDe.Range("H7").Value = "2017-03-02 08:30:09"

# --- Overview sheet: the "Latest HO Xliff Generate Date" column (G) for the rows that
# reference the same generation timestamp moves forward too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-03-02 08:30:09"
$wsOverview.Range("G5").Value = "2017-03-02 08:30:09"
$wsOverview.Range("G6").Value = "2017-03-02 08:30:09"
$wsOverview.Range("G7").Value = "2017-03-02 08:30:09"
